$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.955.01'
$ws.Range('E2').Value = '  -2.17%  '
$ws.Range('D3').Value = '2.421.78'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '550.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.08%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.496'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.85%  '
$ws.Range('D9').Value = '2.417.36'
$ws.Range('E9').Value = '  -2.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.82%  '
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.329'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.83%  '
$ws.Range('D14').Value = '2.863.24'
$ws.Range('E14').Value = '  -2.76%  '
$ws.Range('D15').Value = '67.751.17'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('E16').Value = '  -5.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '22.95'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.66%  '
$ws.Range('D18').Value = '2.406.76'
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.94%  '
$ws.Range('E22').Value = '  -3.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('E24').Value = '  -3.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.00%  '
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.56'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('E29').Value = '  -8.06%  '
$ws.Range('E30').Value = '  -8.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.02%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '427.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('E34').Value = '  -7.11%  '
$ws.Range('E35').Value = '  -6.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '156.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  -4.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('E41').Value = '  -4.98%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.32'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.39%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '37.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('E44').Value = '  -9.40%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '129.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.84%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.97'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.26%  '
$ws.Range('E48').Value = '  -3.97%  '
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.473'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.550'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.71%  '
